$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.131.44"
$ws.Range("E2").Value = "  -3.08%  "
$ws.Range("D3").Value = "1.605.40"
$ws.Range("E3").Value = "  -2.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.000"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "302.29"
$ws.Range("E6").Value = "  -2.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3790"
$ws.Range("E7").Value = "  -2.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3665"
$ws.Range("E8").Value = "  -3.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.25"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.274"
$ws.Range("E10").Value = "  -5.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08154"
$ws.Range("E11").Value = "  -3.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.91"
$ws.Range("E13").Value = "  -4.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.627"
$ws.Range("E14").Value = "  -6.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001263"
$ws.Range("E15").Value = "  -3.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.420"
$ws.Range("E16").Value = "  -7.98%  "
$ws.Range("D17").Value = "1.602.07"
$ws.Range("E17").Value = "  -3.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.27"
$ws.Range("E18").Value = "  -2.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06875"
$ws.Range("E19").Value = "  -1.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.33"
$ws.Range("E20").Value = "  -6.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.611"
$ws.Range("E21").Value = "  -5.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.5552"
$ws.Range("E22").Value = "  -6.08%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.09"
$ws.Range("E24").Value = "  -4.46%  "
$ws.Range("D25").Value = "23.115.27"
$ws.Range("E25").Value = "  -3.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.357"
$ws.Range("E26").Value = "  -2.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.819"
$ws.Range("E27").Value = "  -5.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.19"
$ws.Range("E28").Value = "  -3.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "149.95"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.275"
$ws.Range("E30").Value = "  -2.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "134.02"
$ws.Range("E31").Value = "  -2.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.386"
$ws.Range("E32").Value = "  -3.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.907"
$ws.Range("E33").Value = "  -11.44%  "
$ws.Range("D34").Value = "1.783.73"
$ws.Range("E34").Value = "  -2.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9627"
$ws.Range("E35").Value = "  -4.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.07747"
$ws.Range("E36").Value = "  -5.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.43"
$ws.Range("E37").Value = "  -3.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.332"
$ws.Range("E38").Value = "  -4.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02741"
$ws.Range("E39").Value = "  -5.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2557"
$ws.Range("E40").Value = "  -4.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.08903"
$ws.Range("E41").Value = "  -2.56%  "
$ws.Range("E42").Value = "  -3.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7103"
$ws.Range("E43").Value = "  -6.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.66"
$ws.Range("E44").Value = "  -6.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.41"
$ws.Range("E45").Value = "  -6.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6649"
$ws.Range("E46").Value = "  -4.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.325"
$ws.Range("E47").Value = "  -5.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9988"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.006"
$ws.Range("E49").Value = "  -2.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.85"
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.251"
$ws.Range("E51").Value = "  +1.85%  "
